$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "42.796.09"
$ws.Range("E2").Value = "  -0.02%  "
Set-TextValue $ws.Range("D3") "2.313.94"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "301.25"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.04%  "
Set-TextValue $ws.Range("D9") "0.491"
$ws.Range("E9").Value = "  -1.19%  "
Set-TextValue $ws.Range("D10") "34.09"
$ws.Range("E10").Value = "  -3.05%  "
Set-TextValue $ws.Range("D11") "18.95"
$ws.Range("E11").Value = "  +1.55%  "
Set-TextValue $ws.Range("D12") "0.0782"
$ws.Range("E13").Value = "  +0.49%  "
Set-TextValue $ws.Range("D14") "6.71"
$ws.Range("E14").Value = "  -1.89%  "
Set-TextValue $ws.Range("D15") "2.673.06"
$ws.Range("E15").Value = "  +0.42%  "
Set-TextValue $ws.Range("D16") "2.248.99"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("E17").Value = "  +0.93%  "
Set-TextValue $ws.Range("D18") "42.731.56"
$ws.Range("E18").Value = "  -0.01%  "
Set-TextValue $ws.Range("D19") "12.13"
$ws.Range("E19").Value = "  -4.43%  "
$ws.Range("E20").Value = "  +1.72%  "
Set-TextValue $ws.Range("D21") "0.0₃0890"
$ws.Range("E21").Value = "  -0.43%  "
Set-TextValue $ws.Range("D22") "67.71"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E23").Value = "  +6.27%  "
Set-TextValue $ws.Range("D24") "235.00"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -0.15%  "
Set-TextValue $ws.Range("D27") "24.31"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("E28").Value = "  +14.56%  "
Set-TextValue $ws.Range("D29") "9.11"
$ws.Range("E29").Value = "  +0.88%  "
Set-TextValue $ws.Range("D30") "32.19"
$ws.Range("E30").Value = "  -2.64%  "
Set-TextValue $ws.Range("D31") "148.50"
$ws.Range("E31").Value = "  -10.56%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  +0.49%  "
Set-TextValue $ws.Range("D34") "17.64"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("E35").Value = "  +0.15%  "
Set-TextValue $ws.Range("D36") "0.0699"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  -1.29%  "
Set-TextValue $ws.Range("D38") "1.78"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  -0.95%  "
Set-TextValue $ws.Range("D42") "21.80"
$ws.Range("E42").Value = "  +20.09%  "
Set-TextValue $ws.Range("D43") "1.917.32"
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("E46").Value = "  -1.55%  "
Set-TextValue $ws.Range("D47") "2.74"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("E48").Value = "  +1.56%  "
Set-TextValue $ws.Range("D49") "2.542.20"
$ws.Range("E49").Value = "  +0.48%  "
Set-TextValue $ws.Range("D50") "53.24"
$ws.Range("E50").Value = "  -0.46%  "
Set-TextValue $ws.Range("D51") "72.18"
$ws.Range("E51").Value = "  +1.59%  "
